$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).ColumnWidth = 59.996651785714285
